$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L1").Value = "Elemento1"
$ws.Range("C5").Value = "Costanilla de san Andres, 22 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C6").Value = "Plaza  de San Andres, 2 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C10").Value = "Calle del Aguila, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C15").Value = "Calle de Bailen, 4 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C16").Value = "Calle de Bailen, 6 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C19").Value = "Calle de Barbara de Braganza, 3 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C26").Value = "Plaza  de la Encarnacion, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C27").Value = "Calle de Alcala, 43 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C30").Value = "Plaza  de San Andres, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C31").Value = "Plaza  de Lavapies, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C36").Value = "Calle de Alcala, 25 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C40").Value = "Cra. de San Jeronimo, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C43").Value = "Calle de Tetuan, 23 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C47").Value = "Calle de Gran Via, 17 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C54").Value = "Calle de la Concepcion Jeronima, 15 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C55").Value = "Plaza  de San Nicolas, 6 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C64").Value = "Plaza  de la Armeria, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C65").Value = "Calle del Leon, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C70").Value = "Calle de Bailen, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C74").Value = "Plaza  de Santa Barbara, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C77").Value = "Calle de Ruiz de Alarcon, 23 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C85").Value = "Calle de Barcelo, 2 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C86").Value = "Calle del Meson de Paredes, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C89").Value = "Plaza  de San Martin, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C92").Value = "Plaza  de San Martin, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C95").Value = "Calle de Bailen, 17 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C104").Value = "Calle de Bailen, 7 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C105").Value = "Calle del Principe, 25 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C108").Value = "Calle de Mendez Muñez, 8 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C109").Value = "Plaza  de la Villa de Paris, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C114").Value = "Plaza  de la Villa de Paris, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C115").Value = "Calle de Alcala, 5 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C120").Value = "Plaza  del Cordon, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C133").Value = "Plaza  de Canovas del Castillo, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C137").Value = "Calle de Bailen, 4 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C144").Value = "Calle de Colon, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C147").Value = "Plaza  de la Villa de Paris, 1 ,Madrid, Comunidad de Madrid, España"
$ws.Range("C152").Value = "Paseo de Fernan Nuñez, 4 ,Madrid, Comunidad de Madrid, España"
